$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2222222222222222
$ws.Range("C2").Value = 0.4444444444444444
$ws.Range("J2").Value = 0.05555555555555555
$ws.Range("P2").Value = 0.2222222222222222
$ws.Range("S2").Value = 0.05555555555555555
$ws.Range("J3").Value = 0.125
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.375
$ws.Range("J4").Value = 0.25
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.07142857142857142
$ws.Range("D6").Value = 0.07142857142857142
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("Q6").Value = 0.2857142857142857
$ws.Range("S6").Value = 0.2857142857142857
$ws.Range("F7").Value = 0.125
$ws.Range("J7").Value = 0.25
$ws.Range("R7").Value = 0.125
$ws.Range("S7").Value = 0.5
$ws.Range("B8").Value = 0.09677419354838709
$ws.Range("D8").Value = 0.06451612903225806
$ws.Range("J8").Value = 0.09677419354838709
$ws.Range("O8").Value = 0.06451612903225806
$ws.Range("Q8").Value = 0.1935483870967742
$ws.Range("R8").Value = 0.1290322580645161
$ws.Range("S8").Value = 0.3548387096774194
$ws.Range("D9").Value = 0.125
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.125
$ws.Range("Q9").Value = 0.125
$ws.Range("R9").Value = 0.25
$ws.Range("S9").Value = 0.25
$ws.Range("B10").Value = 0.1075268817204301
$ws.Range("D10").Value = 0.01075268817204301
$ws.Range("F10").Value = 0.1075268817204301
$ws.Range("J10").Value = 0.09677419354838709
$ws.Range("O10").Value = 0.02150537634408602
$ws.Range("Q10").Value = 0.2365591397849462
$ws.Range("R10").Value = 0.1290322580645161
$ws.Range("S10").Value = 0.2903225806451613
$ws.Range("G11").Value = 0.2222222222222222
$ws.Range("J11").Value = 0.05555555555555555
$ws.Range("K11").Value = 0.2777777777777778
$ws.Range("L11").Value = 0.3888888888888889
$ws.Range("S11").Value = 0.05555555555555555
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2857142857142857
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.04166666666666666
$ws.Range("J15").Value = 0.2916666666666667
$ws.Range("K15").Value = 0.125
$ws.Range("O15").Value = 0.25
$ws.Range("S15").Value = 0.125
$ws.Range("H16").Value = 0.1
$ws.Range("J16").Value = 0.6
$ws.Range("K16").Value = 0.1
$ws.Range("S16").Value = 0.2
$ws.Range("H17").Value = 0.1818181818181818
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.4848484848484849
$ws.Range("K17").Value = 0.0303030303030303
$ws.Range("O17").Value = 0.0303030303030303
$ws.Range("S17").Value = 0.1818181818181818
$ws.Range("F18").Value = 0.05555555555555555
$ws.Range("H18").Value = 0.3333333333333333
$ws.Range("I18").Value = 0.05555555555555555
$ws.Range("J18").Value = 0.2777777777777778
$ws.Range("K18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.2222222222222222
$ws.Range("F19").Value = 0.0136986301369863
$ws.Range("H19").Value = 0.1780821917808219
$ws.Range("I19").Value = 0.0410958904109589
$ws.Range("J19").Value = 0.4794520547945205
$ws.Range("K19").Value = 0.0958904109589041
$ws.Range("O19").Value = 0.1232876712328767
$ws.Range("S19").Value = 0.0684931506849315
